# Apply updates to the "two-digit number divided by one-digit number" worksheet.
# 1) Update the date line.
# 2) Update each division problem in the table.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-21 Friday", "2025-03-22 Saturday"),
    @("35÷5=", "92÷3="),
    @("99÷7=", "98÷7="),
    @("35÷3=", "61÷2="),
    @("16÷3=", "88÷4="),
    @("98÷5=", "12÷4="),
    @("15÷3=", "22÷3="),
    @("65÷5=", "41÷6="),
    @("68÷9=", "71÷4="),
    @("94÷3=", "61÷4="),
    @("97÷3=", "79÷9="),
    @("12÷7=", "40÷4="),
    @("19÷2=", "65÷7="),
    @("48÷4=", "73÷7="),
    @("46÷7=", "87÷4="),
    @("79÷5=", "74÷5="),
    @("42÷5=", "43÷2="),
    @("45÷9=", "28÷9="),
    @("95÷7=", "35÷7="),
    @("29÷7=", "50÷9="),
    @("28÷7=", "24÷2="),
    @("11÷7=", "21÷9="),
    @("58÷9=", "33÷9="),
    @("31÷2=", "99÷6="),
    @("53÷6=", "23÷8="),
    @("40÷7=", "66÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
